$d = $word.ActiveDocument

# 1. Locate the "Iterator" Heading-1 paragraph (there are other paragraphs
#    elsewhere that merely start with "Iterator", e.g. the "Cmsor2" ones,
#    so match on the exact paragraph text + style to stay precise) and
#    append " (Behavioral pattern)" to it, mirroring the existing
#    "Singleton (Creational pattern)" / "Visitor (Behavioral pattern)"
#    headings already in the document.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Iterator`r" -and $p.Style.NameLocal -eq "Heading 1") {
        $r = $d.Range($p.Range.Start, $p.Range.End - 1)
        $r.InsertAfter(" (Behavioral")
        $r = $d.Range($r.End, $r.End)
        $r.InsertAfter(" pattern")
        $r = $d.Range($r.End, $r.End)
        $r.InsertAfter(")")
        break
    }
}

# 2. Insert "biztosít " right after "Szekvenciális hozzáférést " in the
#    bullet paragraph describing the Iterator pattern.
$d.Content.Find.Execute("Szekvenciális hozzáférést ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Szekvenciális hozzáférést biztosít ", 2)
